# Update "people" column (B2:B10) with new crossings-per-year data and
# restyle those cells (centered, Helvetica Neue 15pt, thousands separator)
# to match the new look used for the updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 153706850
    3  = 159061181
    4  = 166151041
    5  = 173260603
    6  = 181281186
    7  = 185155513
    8  = 187965778
    9  = 192913686
    10 = 188228921
}

# Set the new values for every row first.
foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value2 = $newValues[$row]
}

# Build the new style once on B2, then propagate it to the rest of the
# column via a format-only paste so only a single new style entry is
# created (instead of one per cell).
$first = $ws.Range("B2")
$first.Font.Name = "Helvetica Neue"
$first.Font.Size = 15
$first.NumberFormat = "#,##0"
$first.HorizontalAlignment = -4108

$first.Copy()
foreach ($row in 3, 4, 5, 6, 7, 8, 9, 10) {
    $ws.Range("B$row").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("B12").Select() | Out-Null
